# "udpate salaries week 3"
# Fill in week-3 salary/date/task data on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B3 (Date) -------------------------------------------------------
# The cell is pre-formatted as a date (m/d/yyyy). The new value
# "03.10.2019" must land as literal TEXT (not get reinterpreted as an
# actual date serial number), while keeping the cell's existing style.
# Entering it through a helper formula cell (so it is already a string,
# not raw user input) and pasting only the value across achieves this
# without creating a stray new cell style.
$helper = $ws.Range("Z1")
$helper.Formula = "=""03.10.2019"""
$helper.Copy()
$ws.Range("B3").PasteSpecial(-4163)
$helper.Clear()

# --- B4 (Team Name) ---------------------------------------------------
$ws.Range("B4").Value = "MSR Voice Input"

# --- B5 (Total Number of Team Members) --------------------------------
$ws.Range("B5").Value = 5

# --- B8:B12 (Member 1-5 salaries) -------------------------------------
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 100

# (B14 total-salary-distributed and B15 total-salary-remaining are
# formulas already on the sheet; they recalc automatically to 500 / 0.)

# --- Tasks completed this week / Tasks to complete next week ---------
$ws.Range("A19").Value = "Forming interview questions"
$ws.Range("A20").Value = "Interview with Sophie"
$ws.Range("B19").Value = "Affinity clustering"
$ws.Range("B20").Value = "Needfinding"
$ws.Range("B21").Value = "Preparing for in-person interview with Sophie"

# --- Cosmetic view / row-height tweaks captured in the authored diff -
$ws.Rows.Item(1).RowHeight = 91
$ws.Rows.Item(18).RowHeight = 41
$excel.ActiveWindow.Zoom = 133
$ws.Range("B22").Select() | Out-Null
